$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.425.94'
$ws.Range('E2').Value = '  -0.50%  '

$ws.Range('D3').Value = '2.616.21'
$ws.Range('E3').Value = '  +0.03%  '

$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.56%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.552'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.49%  '

$ws.Range('D9').Value = '2.614.56'
$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('E10').Value = '  -2.23%  '

$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('E12').Value = '  -0.78%  '

$ws.Range('E13').Value = '  -2.00%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.46%  '

$ws.Range('D15').Value = '3.091.07'
$ws.Range('E15').Value = '  +0.43%  '

$ws.Range('E16').Value = '  -3.27%  '

$ws.Range('D17').Value = '67.299.75'
$ws.Range('E17').Value = '  -0.50%  '

$ws.Range('D18').Value = '2.616.30'
$ws.Range('E18').Value = '  +0.24%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '367.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.55%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.32%  '

$ws.Range('E21').Value = '  -3.86%  '

$ws.Range('E22').Value = '  -0.32%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.71%  '

$ws.Range('E25').Value = '  -0.03%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.06'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '66.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.91%  '

$ws.Range('E28').Value = '  +0.39%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '582.40'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.91%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.17%  '

$ws.Range('E31').Value = '  -3.30%  '

$ws.Range('E32').Value = '  -3.66%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.26%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.73%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.126'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.17%  '

$ws.Range('E37').Value = '  -1.53%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '156.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.51%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.03'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.46%  '

$ws.Range('E40').Value = '  +2.77%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.366'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.06%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.78%  '

$ws.Range('E43').Value = '  -0.84%  '

$ws.Range('E44').Value = '  +2.49%  '

$ws.Range('E45').Value = '  -0.05%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '155.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('D47').Value = '0.0₆0293'
$ws.Range('E47').Value = '  +2.33%  '

$ws.Range('E48').Value = '  -0.76%  '

$ws.Range('E49').Value = '  -1.09%  '

$ws.Range('E50').Value = '  +0.58%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.86%  '
